$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2 = @{ E = 3; G = 4.619088000000001; H = 13.857264; K = 2; M = 17.945525; N = 35.89105; O = 0.3147738875783; P = 0.2472168478181395; Q = 82.89195918120001; R = 497.3517550872; S = 0.3147738875783; T = 0.2472168478181395 }
    3 = @{ E = 3; G = 4.619088000000001; H = 13.857264; K = 3; M = 3.914977; N = 11.744931; O = 0.06867074270993077; P = 0.0808988541617353; Q = 18.083623280976; R = 162.752609528784; S = 0.06867074270993077; T = 0.0808988541617353 }
    4 = @{ E = 3; G = 4.619088000000001; H = 13.857264; K = 3; M = 7.595080333333333; N = 22.785241; O = 0.133221678551774; P = 0.1569442927079769; Q = 35.082344426736; R = 315.741099840624; S = 0.133221678551774; T = 0.1569442927079769 }
    5 = @{ E = 3; G = 4.619088000000001; H = 13.857264; K = 3; M = 8.835736333333331; N = 26.507209; O = 0.1549834332102386; P = 0.18258113522554; Q = 40.813043668464; R = 367.317393016176; S = 0.1549834332102386; T = 0.18258113522554 }
    6 = @{ E = 3; G = 4.619088000000001; H = 13.857264; K = 3; M = 10.81295533333333; N = 32.438866; O = 0.1896648878471846; P = 0.2234382721964117; Q = 49.94599222473601; R = 449.5139300226241; S = 0.1896648878471846; T = 0.2234382721964117 }
    7 = @{ E = 3; G = 4.619088000000001; H = 13.857264; K = 2; M = 7.90657; N = 15.81314; O = 0.1386853701025721; P = 0.1089205978901965; Q = 36.52114260816001; R = 219.12685564896; S = 0.1386853701025721; T = 0.1089205978901965 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}
